$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D and E to remain plain text so numeric-looking values
# (prices, percentages) are not auto-converted to Number cells.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "61.117.76"
$ws.Range("E2").Value = "  -4.32%  "

# Row 3
$ws.Range("D3").Value = "2.453.82"
$ws.Range("E3").Value = "  -6.88%  "

# Row 4
$ws.Range("E4").Value = "  +0.07%  "

# Row 5
$ws.Range("D5").Value = "546.01"
$ws.Range("E5").Value = "  -5.77%  "

# Row 6
$ws.Range("D6").Value = "145.86"
$ws.Range("E6").Value = "  -7.05%  "

# Row 7
$ws.Range("E7").Value = "  +0.05%  "

# Row 8
$ws.Range("D8").Value = "0.584"
$ws.Range("E8").Value = "  -7.68%  "

# Row 9
$ws.Range("D9").Value = "2.451.85"
$ws.Range("E9").Value = "  -6.92%  "

# Row 10
$ws.Range("E10").Value = "  -10.73%  "

# Row 11
$ws.Range("B11").Value = "Toncoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D11").Value = "5.43"
$ws.Range("E11").Value = "  -6.92%  "

# Row 12
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "0.154"
$ws.Range("E12").Value = "  -1.97%  "

# Row 13
$ws.Range("D13").Value = "0.351"
$ws.Range("E13").Value = "  -8.94%  "

# Row 14
$ws.Range("D14").Value = "25.93"
$ws.Range("E14").Value = "  -9.86%  "

# Row 15
$ws.Range("D15").Value = "2.897.92"
$ws.Range("E15").Value = "  -7.05%  "

# Row 16
$ws.Range("D16").Value = "0.0000166"
$ws.Range("E16").Value = "  -10.42%  "

# Row 17
$ws.Range("D17").Value = "61.098.81"
$ws.Range("E17").Value = "  -4.27%  "

# Row 18
$ws.Range("D18").Value = "2.448.59"
$ws.Range("E18").Value = "  -7.72%  "

# Row 19
$ws.Range("D19").Value = "11.06"
$ws.Range("E19").Value = "  -9.38%  "

# Row 20
$ws.Range("E20").Value = "  -9.52%  "

# Row 21
$ws.Range("D21").Value = "4.15"
$ws.Range("E21").Value = "  -8.44%  "

# Row 22
$ws.Range("D22").Value = "316.67"
$ws.Range("E22").Value = "  -8.04%  "

# Row 24
$ws.Range("E24").Value = "  -2.49%  "

# Row 25
$ws.Range("D25").Value = "63.71"
$ws.Range("E25").Value = "  -6.85%  "

# Row 26
$ws.Range("D26").Value = "2.579.96"
$ws.Range("E26").Value = "  -6.87%  "

# Row 27
$ws.Range("D27").Value = "550.14"
$ws.Range("E27").Value = "  -5.46%  "

# Row 28
$ws.Range("D28").Value = "0.0₃0957"
$ws.Range("E28").Value = "  -15.17%  "

# Row 29
$ws.Range("E29").Value = "  -2.71%  "

# Row 30
$ws.Range("E30").Value = "  -12.27%  "

# Row 31
$ws.Range("D31").Value = "8.18"
$ws.Range("E31").Value = "  -11.71%  "

# Row 32
$ws.Range("D32").Value = "7.53"
$ws.Range("E32").Value = "  -9.46%  "

# Row 33
$ws.Range("D33").Value = "0.146"
$ws.Range("E33").Value = "  -9.00%  "

# Row 34
$ws.Range("E34").Value = "  -9.31%  "

# Row 35
$ws.Range("D35").Value = "1.58"
$ws.Range("E35").Value = "  -9.83%  "

# Row 36
$ws.Range("D36").Value = "5.83"
$ws.Range("E36").Value = "  -12.87%  "

# Row 37
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  -0.01%  "

# Row 38
$ws.Range("D38").Value = "4.79"
$ws.Range("E38").Value = "  -12.76%  "

# Row 39
$ws.Range("E39").Value = "  -7.01%  "

# Row 40
$ws.Range("D40").Value = "18.38"
$ws.Range("E40").Value = "  -7.29%  "

# Row 41
$ws.Range("D41").Value = "1.76"
$ws.Range("E41").Value = "  -8.27%  "

# Row 42
$ws.Range("D42").Value = "140.93"
$ws.Range("E42").Value = "  -8.62%  "

# Row 43
$ws.Range("E43").Value = "  +0.04%  "

# Row 44
$ws.Range("D44").Value = "40.39"
$ws.Range("E44").Value = "  -4.21%  "

# Row 45
$ws.Range("D45").Value = "2.29"
$ws.Range("E45").Value = "  -11.05%  "

# Row 46
$ws.Range("D46").Value = "145.93"
$ws.Range("E46").Value = "  -10.70%  "

# Row 47
$ws.Range("D47").Value = "3.58"
$ws.Range("E47").Value = "  -9.06%  "

# Row 48
$ws.Range("D48").Value = "21.34"
$ws.Range("E48").Value = "  -12.10%  "

# Row 49
$ws.Range("D49").Value = "0.0533"
$ws.Range("E49").Value = "  -9.62%  "

# Row 50
$ws.Range("D50").Value = "0.586"
$ws.Range("E50").Value = "  -7.58%  "

# Row 51
$ws.Range("D51").Value = "0.0928"
